$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Sheet" to "Sheet1"
$ws.Name = "Sheet1"

# Update existing row 2 values
$ws.Range("A2").Value = "Item1"
$ws.Range("B2").Value = 100

# Add new row 3
$ws.Range("A3").Value = "Item2"
$ws.Range("B3").Value = 200
